$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove now-unused placeholder rows 18:20 (table shrinks from 19 to 16 data rows) ---
$ws.Range("A18:N20").EntireRow.Delete() | Out-Null

# --- New text cells, written in the exact order the original author entered them ---
# --- (this keeps freshly-appended shared-string indices aligned with the target file) ---
$t1 = @'
Network Ports
'@
$ws.Range("B12").Value = $t1
$t2 = @'
Server and client not able to communicate with each other
'@
$ws.Range("C12").Value = $t2
$t3 = @'
1. Server and client not able to communicate with each other
2. System is not usable at all
'@
$ws.Range("D12").Value = $t3
$t4 = @'
1. Network ports are not available is used by another process
2. Firewalls is blocking the port 
'@
$ws.Range("F12").Value = $t4
$t5 = @'
1. System is not usable, User interface will show the network
2. No REST API is called on the server side
'@
$ws.Range("H12").Value = $t5
$t6 = @'
1. Request IT to open the desired ports
2. Make the PORT configurable via the environment variables
3. Ask IT to whitelist the PORT in the firewall

'@
$ws.Range("K12").Value = $t6
$t7 = @'
Secured Communication
'@
$ws.Range("B13").Value = $t7
$t8 = @'
Communication between the server and client is not secured
'@
$ws.Range("C13").Value = $t8
$t9 = @'
1. The data can easily be hijacked or tempared over the wire.
2. Loss of sensitive organisation data
'@
$ws.Range("D13").Value = $t9
$t10 = @'
1. The system is running on the public internet
2. System is not using the HTTPS protocol for communication
'@
$ws.Range("F13").Value = $t10
$t11 = @'
1. The Headers of the REST calls are using HTTP headers
2. One can access the system over the public internet
'@
$ws.Range("H13").Value = $t11
$t12 = @'
1. Use the HTTPS protocl for communication between server and client
2. Install the system behind the VPN so only authorized people can access it
'@
$ws.Range("K13").Value = $t12
$t13 = @'
Export Information Security
'@
$ws.Range("B14").Value = $t13
$t14 = @'
The data exported out of the app is not secured and encrypted
'@
$ws.Range("C14").Value = $t14
$t15 = @'
1. The exported data can easily be read or modified by any user.
2. The exported data can easily be mis-used
'@
$ws.Range("D14").Value = $t15
$t16 = @'
1. The sensitive data in the exported is not encrypted.
2. Exported files can be opened in any freely available editors
3. Files are not password protected
'@
$ws.Range("F14").Value = $t16
$t17 = @'
1. Files can easily be exported out, read in any freely avaibale editor
'@
$ws.Range("H14").Value = $t17
$t18 = @'
1. Encrypt the sensitive data while exporting the file,
2. Make the file passowrd protected and only authorized person can open the files.
3. Change the extension of the file so that not every editor can open and develop a small editor for the file in the application
'@
$ws.Range("K14").Value = $t18
$t19 = @'
User Interface Dashboard
'@
$ws.Range("B15").Value = $t19
$t20 = @'
Graphs are not easily understood by the operator. The Graphs represent the machine vitals.
'@
$ws.Range("C15").Value = $t20
$t21 = @'
1. Operator not able deduce the correct informaton from the graphs, as graphs representation is good for the operator
'@
$ws.Range("D15").Value = $t21
$t22 = @'
1. Operator will complain about the usage of the app.
2. Operator might take wrong decisions
'@
$ws.Range("H15").Value = $t22
$t23 = @'
1. Graphs representation is not known to operator
'@
$ws.Range("F15").Value = $t23
$t24 = @'
1. Provide the different options to visualize the machine vitals like in tables
'@
$ws.Range("K15").Value = $t24
$t25 = @'
1. Not able to visualize any of the machine vitals in the dashboard
2. Prediction computation for the machine is impossible
'@
$ws.Range("D16").Value = $t25
$t26 = @'
1.Machine Id entered during the machine addition is incorrect
'@
$ws.Range("F16").Value = $t26
$t27 = @'
New Machine Addition
'@
$ws.Range("B16").Value = $t27
$t28 = @'
New machine is not properly added to the system
'@
$ws.Range("C16").Value = $t28
$t29 = @'
1. Data is present is the telemetry file but can not be visualzied in the dashboard
2. Getting no prediction classification for the machine
'@
$ws.Range("H16").Value = $t29
$t30 = @'
Deployment
'@
$ws.Range("B17").Value = $t30
$t31 = @'
System is not getting deployed in new environment
'@
$ws.Range("C17").Value = $t31
$t32 = @'
1. System is giving errors during the deployment as required packages are not getting installed
'@
$ws.Range("H17").Value = $t32
$t33 = @'
1. System is not deployed and hence can not be started in the target environment
'@
$ws.Range("D17").Value = $t33
$t34 = @'
1. No internet is available to download the desired packages
'@
$ws.Range("F17").Value = $t34
$t35 = @'
1. Ask the IT guys to make the internet available on the system before starting the deployment.
2. Package the desired modules as part of the software package
'@
$ws.Range("K17").Value = $t35
$t36 = @'
1. Provide the operator a proper machine manually
2. Create a REST API that will ping the all the machine connected to system to get the machine id
3. Create a method in the server that will read the telemetry file for all unique machine id's available in the file
4. Train the operator for using the new software
'@
$ws.Range("K16").Value = $t36

# --- Numeric severity/occurrence/probability cells ---
$ws.Range("E12").Value = 10
$ws.Range("G12").Value = 2
$ws.Range("I12").Value = 1
$ws.Range("E13").Value = 10
$ws.Range("G13").Value = 8
$ws.Range("I13").Value = 3
$ws.Range("E14").Value = 10
$ws.Range("G14").Value = 10
$ws.Range("I14").Value = 1
$ws.Range("E15").Value = 8
$ws.Range("G15").Value = 4
$ws.Range("I15").Value = 1
$ws.Range("E16").Value = 10
$ws.Range("G16").Value = 5
$ws.Range("I16").Value = 1
$ws.Range("E17").Value = 7
$ws.Range("G17").Value = 3
$ws.Range("I17").Value = 2

# --- J column: shared RPN formula (E*G*I) across J2:J17, anchored at J2 ---
$ws.Range("J2:J17").Formula = "=(E2*G2*I2)"

# --- Row heights for the newly populated rows ---
$ws.Rows.Item(12).RowHeight = 115
$ws.Rows.Item(13).RowHeight = 102
$ws.Rows.Item(14).RowHeight = 170
$ws.Rows.Item(15).RowHeight = 85
$ws.Rows.Item(16).RowHeight = 221
$ws.Rows.Item(17).RowHeight = 119

# --- Restore selection near the newly added content ---
$ws.Range("K19").Select() | Out-Null

